$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1760.6923
$ws.Range("J17").Value = 1760.6923
$ws.Range("L17").Value = 5282.0769
$ws.Range("N17").Value = -5618.0769
$ws.Range("H64").Value = 6987.778
$ws.Range("I64").Value = 5198
$ws.Range("K64").Value = 5198
$ws.Range("M64").Value = -4950
$ws.Range("H67").Value = 6987.778
$ws.Range("I67").Value = 5198
$ws.Range("K67").Value = 5198
$ws.Range("M67").Value = -4340
$ws.Range("H96").Value = 7937325.5
$ws.Range("I96").Value = 17857370
$ws.Range("J96").Value = 1289
$ws.Range("K96").Value = 53572110
$ws.Range("L96").Value = 3867
$ws.Range("M96").Value = -53570737
$ws.Range("N96").Value = -6613
$ws.Range("H98").Value = 1511.5526
$ws.Range("I98").Value = 1570.3125
$ws.Range("J98").Value = 1198.1666
$ws.Range("K98").Value = 1570.3125
$ws.Range("L98").Value = 1198.1666
$ws.Range("M98").Value = -72.3125
$ws.Range("N98").Value = -4194.1666
$ws.Range("H113").Value = 5314.6
$ws.Range("I113").Value = 5314.6
$ws.Range("K113").Value = 5314.6
$ws.Range("M113").Value = -2060.6
$ws.Range("H122").Value = 1511.5526
$ws.Range("I122").Value = 1570.3125
$ws.Range("J122").Value = 1198.1666
$ws.Range("K122").Value = 4710.9375
$ws.Range("L122").Value = 3594.4998
$ws.Range("M122").Value = -2260.9375
$ws.Range("N122").Value = -8494.4998
$ws.Range("H132").Value = 2709806.2
$ws.Range("I132").Value = 2761881.2
$ws.Range("K132").Value = 8285643.600000001
$ws.Range("M132").Value = -8283113.600000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2371.6
$ws.Range("J45").Value = 4964.3335
$ws.Range("L45").Value = 4964.3335
$ws.Range("N45").Value = -5718.3335
$ws.Range("H61").Value = 6171.84
$ws.Range("I61").Value = 923.5833
$ws.Range("K61").Value = 923.5833
$ws.Range("M61").Value = -711.5833
$ws.Range("H74").Value = 122928.42
$ws.Range("I74").Value = 140184.81
$ws.Range("J74").Value = 16924.857
$ws.Range("K74").Value = 140184.81
$ws.Range("L74").Value = 16924.857
$ws.Range("M74").Value = -139310.81
$ws.Range("N74").Value = -18672.857
$ws.Range("H77").Value = 122928.42
$ws.Range("I77").Value = 140184.81
$ws.Range("J77").Value = 16924.857
$ws.Range("K77").Value = 700924.05
$ws.Range("L77").Value = 84624.285
$ws.Range("M77").Value = -696556.05
$ws.Range("N77").Value = -93360.285
$ws.Range("H102").Value = 2735.6428
$ws.Range("I102").Value = 2791.3845
$ws.Range("K102").Value = 2791.3845
$ws.Range("M102").Value = -1169.3845
$ws.Range("H122").Value = 1502.9688
$ws.Range("I122").Value = 1502.9688
$ws.Range("K122").Value = 4508.9064
$ws.Range("M122").Value = -2058.9064
$ws.Range("H136").Value = 6171.84
$ws.Range("I136").Value = 923.5833
$ws.Range("K136").Value = 2770.7499
$ws.Range("M136").Value = -220.7498999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 9999.5
$ws.Range("J6").Value = 9999.5
$ws.Range("L6").Value = 9999.5
$ws.Range("N6").Value = -10225.5
$ws.Range("H105").Value = 4237
$ws.Range("I105").Value = 3140
$ws.Range("J105").Value = 4895.2
$ws.Range("K105").Value = 3140
$ws.Range("L105").Value = 4895.2
$ws.Range("M105").Value = -1393
$ws.Range("N105").Value = -8389.200000000001
$ws.Range("H123").Value = 99999
$ws.Range("I123").Value = 99999
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 99999
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -95099
$ws.Range("N123").Value = $null  # was -109800
$ws.Range("H134").Value = 1636.9333
$ws.Range("I134").Value = 1202.6052
$ws.Range("K134").Value = 3607.8156
$ws.Range("M134").Value = -1072.8156

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 711.55554
$ws.Range("J22").Value = 794.7143
$ws.Range("L22").Value = 794.7143
$ws.Range("N22").Value = -1494.7143
$ws.Range("H31").Value = 3335683.5
$ws.Range("I31").Value = 4546528.5
$ws.Range("K31").Value = 4546528.5
$ws.Range("M31").Value = -4546233.5
$ws.Range("H34").Value = 3335683.5
$ws.Range("I34").Value = 4546528.5
$ws.Range("K34").Value = 4546528.5
$ws.Range("M34").Value = -4546326.5
$ws.Range("H62").Value = 5999.3335
$ws.Range("I62").Value = 6219.7646
$ws.Range("K62").Value = 6219.7646
$ws.Range("M62").Value = -5595.7646
$ws.Range("H65").Value = 5999.3335
$ws.Range("I65").Value = 6219.7646
$ws.Range("K65").Value = 31098.823
$ws.Range("M65").Value = -27978.823
$ws.Range("H94").Value = 1652.9131
$ws.Range("I94").Value = 1405.4445
$ws.Range("J94").Value = 1812
$ws.Range("K94").Value = 1405.4445
$ws.Range("L94").Value = 1812
$ws.Range("M94").Value = -954.4445000000001
$ws.Range("N94").Value = -2714
$ws.Range("H122").Value = 1485.4546
$ws.Range("I122").Value = 1432
$ws.Range("K122").Value = 4296
$ws.Range("M122").Value = -1846
$ws.Range("H131").Value = 90178.5
$ws.Range("J131").Value = 90178.5
$ws.Range("L131").Value = 90178.5
$ws.Range("N131").Value = -100258.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 765.7917
$ws.Range("I97").Value = 565.8333
$ws.Range("J97").Value = 832.44446
$ws.Range("K97").Value = 1697.4999
$ws.Range("L97").Value = 2497.33338
$ws.Range("M97").Value = -1201.4999
$ws.Range("N97").Value = -3489.33338
$ws.Range("H137").Value = 3127.9285
$ws.Range("J137").Value = 3545.7778
$ws.Range("L137").Value = 10637.3334
$ws.Range("N137").Value = -20837.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 48785.145
$ws.Range("J101").Value = 48785.145
$ws.Range("L101").Value = 48785.145
$ws.Range("N101").Value = -55275.145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1256.7693
$ws.Range("I22").Value = 1104.75
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1104.75
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -809.75
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1256.7693
$ws.Range("I27").Value = 1104.75
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 1104.75
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -997.75
$ws.Range("N27").Value = -1714
$ws.Range("H40").Value = 1982.7142
$ws.Range("I40").Value = 1870.963
$ws.Range("K40").Value = 1870.963
$ws.Range("M40").Value = -1734.963
$ws.Range("H100").Value = 4000
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
$ws.Range("H122").Value = 7064.778
$ws.Range("I122").Value = 5260
$ws.Range("K122").Value = 15780
$ws.Range("M122").Value = -13330

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1525.1111
$ws.Range("I100").Value = 1563
$ws.Range("J100").Value = 1506.1666
$ws.Range("K100").Value = 3126
$ws.Range("L100").Value = 3012.3332
$ws.Range("M100").Value = -2585
$ws.Range("N100").Value = -4094.3332
$ws.Range("H107").Value = 878.5454999999999
$ws.Range("I107").Value = 830.1429000000001
$ws.Range("K107").Value = 2490.4287
$ws.Range("M107").Value = -570.4287000000004
$ws.Range("H132").Value = 1445.7941
$ws.Range("I132").Value = 846.6087
$ws.Range("J132").Value = 2698.6365
$ws.Range("K132").Value = 2539.8261
$ws.Range("L132").Value = 8095.9095
$ws.Range("M132").Value = -9.826100000000224
$ws.Range("N132").Value = -13155.9095
$ws.Range("H137").Value = 104387.86
$ws.Range("J137").Value = 116143
$ws.Range("L137").Value = 116143
$ws.Range("N137").Value = -126343
